$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "1.006") must be
# forced to Text so they stay inline/shared strings like the source data
# (the site renders thousands-grouped prices as text, not numbers).
$numericLookingCells = @(
    @{Cell='D4'; Value='1.006'},
    @{Cell='D5'; Value='217.93'},
    @{Cell='D6'; Value='0.5082'},
    @{Cell='D7'; Value='1.006'},
    @{Cell='D8'; Value='0.2644'},
    @{Cell='D9'; Value='0.06339'},
    @{Cell='D10'; Value='21.46'},
    @{Cell='D11'; Value='0.07360'},
    @{Cell='D13'; Value='4.543'},
    @{Cell='D14'; Value='0.5770'},
    @{Cell='D16'; Value='0.000008500'},
    @{Cell='D17'; Value='64.93'},
    @{Cell='D19'; Value='4.943'},
    @{Cell='D20'; Value='1.006'},
    @{Cell='D21'; Value='10.81'},
    @{Cell='D22'; Value='189.14'},
    @{Cell='D23'; Value='6.192'},
    @{Cell='D24'; Value='1.007'},
    @{Cell='D25'; Value='142.99'},
    @{Cell='D26'; Value='7.660'},
    @{Cell='D27'; Value='0.1174'},
    @{Cell='D28'; Value='15.66'},
    @{Cell='D29'; Value='0.05809'},
    @{Cell='D30'; Value='1.280'},
    @{Cell='D32'; Value='3.500'},
    @{Cell='D33'; Value='3.502'},
    @{Cell='D34'; Value='1.655'},
    @{Cell='D35'; Value='1.006'},
    @{Cell='D36'; Value='0.5982'},
    @{Cell='D38'; Value='2.635'},
    @{Cell='D39'; Value='0.01607'},
    @{Cell='D40'; Value='5.995'},
    @{Cell='D42'; Value='0.8579'},
    @{Cell='D44'; Value='99.71'},
    @{Cell='D46'; Value='0.00000000111'},
    @{Cell='D47'; Value='55.65'},
    @{Cell='D48'; Value='1.004'},
    @{Cell='D49'; Value='8.042'},
    @{Cell='D50'; Value='0.4299'},
    @{Cell='D51'; Value='0.05180'}
)

foreach ($item in $numericLookingCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Remaining cells (plain text, or text containing extra "." separators / % signs)
# can be assigned directly without special handling.
$ws.Range('D2').Value = '26.116.03'
$ws.Range('E2').Value = '  -6.77%  '
$ws.Range('D3').Value = '1.667.19'
$ws.Range('E3').Value = '  -4.33%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('E5').Value = '  -3.75%  '
$ws.Range('E6').Value = '  -12.16%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('E8').Value = '  -2.83%  '
$ws.Range('E9').Value = '  -3.97%  '
$ws.Range('E10').Value = '  -7.39%  '
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('D12').Value = '1.667.94'
$ws.Range('E12').Value = '  -4.51%  '
$ws.Range('E13').Value = '  -3.57%  '
$ws.Range('E14').Value = '  -4.34%  '
$ws.Range('D15').Value = '1.897.75'
$ws.Range('E15').Value = '  -4.19%  '
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('E17').Value = '  -12.87%  '
$ws.Range('D18').Value = '26.187.60'
$ws.Range('E18').Value = '  -6.51%  '
$ws.Range('E19').Value = '  -7.18%  '
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  -4.18%  '
$ws.Range('E22').Value = '  -7.90%  '
$ws.Range('E23').Value = '  -6.57%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  -4.70%  '
$ws.Range('E26').Value = '  -5.50%  '
$ws.Range('E27').Value = '  -4.77%  '
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('E29').Value = '  -5.31%  '
$ws.Range('E30').Value = '  -7.65%  '
$ws.Range('E31').Value = '  -5.23%  '
$ws.Range('E32').Value = '  -6.50%  '
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -3.00%  '
$ws.Range('E36').Value = '  -6.44%  '
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E40').Value = '  -3.28%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.085.28'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').Value = '1.820.41'
$ws.Range('E45').Value = '  -3.80%  '
$ws.Range('E46').Value = '  +3.91%  '
$ws.Range('E47').Value = '  -6.35%  '
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('E49').Value = '  -2.91%  '
$ws.Range('E50').Value = '  -2.74%  '
$ws.Range('E51').Value = '  -3.71%  '
